$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.207052946090698
$ws.Range("B1").Value = 2.798991680145264
$ws.Range("C1").Value = 2.956960678100586
$ws.Range("D1").Value = 2.594427585601807
$ws.Range("E1").Value = 0.8945170640945435
